# Update cryptos list: apply latest price/volume scrape values
# scraped and committed via GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.941.19'
$ws.Range('E2').Value = '  -1.75%  '
$ws.Range('D3').Value = '3.950.61'
$ws.Range('E3').Value = '  -2.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.77'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.29'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('D7').Value = '3.950.83'
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('E8').Value = '  -6.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.738'
$ws.Range('E10').Value = '  -4.90%  '
$ws.Range('E11').Value = '  -6.21%  '
$ws.Range('E12').Value = '  +13.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000316'
$ws.Range('E13').Value = '  -3.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.60'
$ws.Range('E14').Value = '  -5.07%  '
$ws.Range('D15').Value = '4.567.63'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = '3.946.32'
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.60'
$ws.Range('E17').Value = '  -3.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.80'
$ws.Range('E18').Value = '  -2.89%  '
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.17'
$ws.Range('E20').Value = '  -4.22%  '
$ws.Range('D21').Value = '70.802.88'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '425.57'
$ws.Range('E22').Value = '  -4.22%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.58'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '97.01'
$ws.Range('E24').Value = '  -7.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.20'
$ws.Range('E25').Value = '  +4.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.43'
$ws.Range('E26').Value = '  -3.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.34'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.64'
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('E29').Value = '  +15.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.89'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.40'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.74'
$ws.Range('E32').Value = '  +14.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.19'
$ws.Range('E33').Value = '  +17.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.131'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.33'
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '683.23'
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '65.29'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.437'
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('D39').Value = '0.0₃0816'
$ws.Range('E39').Value = '  -5.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.150'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0480'
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('E46').Value = '  -7.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.68'
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.83'
$ws.Range('E48').Value = '  +5.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.31'
$ws.Range('E49').Value = '  -6.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.99'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('E51').Value = '  -1.29%  '
